# Create the new "OutrightHPNRwithoutMaintenance" scenario sheet as the
# last tab in the workbook (this also makes it the active/selected sheet,
# matching Excel's default behaviour when a sheet is inserted via the UI).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "OutrightHPNRwithoutMaintenance"

# Header row
$headers = @(
    "manufacurer",
    "model",
    "Vehicle_Basic_price",
    "road_tax_for_first_year",
    "percentage_cap_residual_value",
    "residual_value_used",
    "additional_terms",
    "additional_mileage",
    "vehicle_profit",
    "maintenance_status",
    "matrix_credit_type",
    "security_deposit",
    "balloon_payment_status",
    "part_exchange_actual",
    "part_exchange_given",
    "less_finance_settlement",
    "order_deposit",
    "finance_deposit",
    "document_fee",
    "sheet_name"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Sample data row
$row2 = @(
    "AUDI",
    "A7",
    50000,
    630,
    30,
    20000,
    22,
    22000,
    1000,
    "NO",
    "A1 Credit",
    100,
    "YES",
    2000,
    1000,
    0,
    0,
    0,
    199.99,
    "HP(N) - HP(N),HP,CP,PCP"
)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}
